$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "Java Notes"/"JavaScript Notes" columns
# (old D -> new F, old E -> new G), making room for new "mp4?" and "ver5 OK?" columns.
$ws.Columns("D:E").Insert()

# Update header row
$ws.Range("A1").Value = "Online?"
$ws.Range("D1").Value = "mp4?"
$ws.Range("E1").Value = "ver5 OK?"

# Update the "2-D Collision of Pucks" entry title and mark new columns as verified
$ws.Range("B2").Value = "2-D Collision of Pucks -- Center of Mass Velocity 2011, Aaron Titus"
$ws.Range("D2").Value = "yes"
$ws.Range("E2").Value = "yes"

# Mark the "Tracker Video Analysis: Remote Control Helicopter" row as verified in new column E
$ws.Range("E48").Value = "yes"
